$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1. Remove the duplicate row (old row 8: "AC-RFD-00-125-000-P1B1-01" /
#    "125 KHz" variant that has no hyperlink - this was a mistaken duplicate
#    of the correct row at row 3).
# -------------------------------------------------------------------------
$ws.Rows("8:8").Delete()

# -------------------------------------------------------------------------
# 2. Remove the stray trailing empty row (old row 13, now row 12 after the
#    deletion above).
# -------------------------------------------------------------------------
$ws.Rows("12:12").Delete()

# -------------------------------------------------------------------------
# 3. Insert a brand new row right after the header and fill it with the
#    new product entry.
# -------------------------------------------------------------------------
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "AC-AKB-00-000-000-P1B1-01"
$ws.Range("C2").Value = "Akbil"
$ws.Range("D2").Value = "Yok"
$ws.Range("E2").Value = "Yok"
$ws.Range("F2").Value = "Genel"
$ws.Range("G2").Value = "Paralel"
$ws.Range("H2").Value = "Vidalı"
$ws.Range("I2").Value = "Buzzerlı"
$ws.Range("J2").Value = "Model-01"
$ws.Range("K2").Value = ""

# Copy the visual formatting of the row below (the previous "first" data
# row) onto the newly inserted row so it looks consistent with the rest of
# the table.
$ws.Range("A3:K3").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values (paste-special formats could have touched values too
# on some hosts) to be safe.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "AC-AKB-00-000-000-P1B1-01"
$ws.Range("C2").Value = "Akbil"
$ws.Range("D2").Value = "Yok"
$ws.Range("E2").Value = "Yok"
$ws.Range("F2").Value = "Genel"
$ws.Range("G2").Value = "Paralel"
$ws.Range("H2").Value = "Vidalı"
$ws.Range("I2").Value = "Buzzerlı"
$ws.Range("J2").Value = "Model-01"
$ws.Range("K2").Value = ""

# -------------------------------------------------------------------------
# 4. Renumber the "Sıra" column (A) sequentially for every data row.
# -------------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# -------------------------------------------------------------------------
# 5. Rebuild the hyperlinks from scratch so every link lands on the right
#    cell after the row shuffling above.
# -------------------------------------------------------------------------
$ws.Range("K2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("K3"), "https://github.com/btk42/AC-AKB-18-000-000-H3B1-01") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K4"), "https://github.com/btk42/AC-RFD-00-125-000-P1B1-01") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K5"), "https://github.com/btk42/AC-RFD-18-125-000-H3B1-01") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K6"), "https://github.com/btk42/AC-RFD-EK-135-LOP-S2B1-01") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K7"), "https://github.com/btk42/AC-RFD-EK-135-COP-S2B1-01") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K8"), "https://github.com/btk42/AC-RFD-EK-135-ANT-H4B0-01") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K9"), "https://github.com/btk42/AC-RFD-00-235-000-H3B1-01") | Out-Null

# -------------------------------------------------------------------------
# 6. Keep the selection in a sane place, mirroring the author's last view.
# -------------------------------------------------------------------------
$ws.Range("E15").Select()

$wb.Save()
